$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.515.51"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.88"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.73"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3662"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07281"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8649"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.65"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.67%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.895.57"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.403"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.545"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06954"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.08%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.003"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "80.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008917"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.42"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.631.45"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.144"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.86"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +5.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.099.95"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.979"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.56"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.87%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.90"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.130"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.27"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.80%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.838"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08865"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.31%  "
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7516"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.15%  "
$ws.Range("B33").Value = "HuobiToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.990"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.545"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.134"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.54%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -1.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05330"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01939"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.795"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5084"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1663"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.560"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.333"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.44"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "105.93"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06497"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4686"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.45%  "
$ws.Range("E49").Value = "  -0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.612"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.83"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.39%  "
